{"js": "// Change summary (per the XML diff):\n//  1. First  \"{{ requestor.preferred_phone_number }}\"  ->  \"{{ requestor.phone_number }}\"\n//  2. Second \"{{ requestor.preferred_phone_number }}\"  ->  \"{{ requestor. phone_number }}\"\n//     (note the stray space after the dot, copied verbatim from the diff)\n//  3. The \"_GoBack\" bookmark is removed from its old spot (right before the\n//     second \"{% if student.grade ... %}\" block) and re-inserted just after\n//     \"requestor.\" inside the second phone-number placeholder above.\n\nconst body = context.document.body;\n\n// Load every paragraph's text up front so we can unambiguously locate the\n// two paragraphs this edit touches.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet firstPhoneParagraph = null;   // \"{{ requestor.preferred_phone_number }} \"\nlet secondPhoneParagraph = null;  // \"Thank you ... {{ requestor.preferred_phone_number }} with ...\"\n\nfor (const p of paragraphs.items) {\n  if (firstPhoneParagraph === null && /^\\{\\{\\s*requestor\\.preferred_phone_number\\b/.test(p.text)) {\n    firstPhoneParagraph = p;\n  }\n  if (p.text.indexOf(\"Thank you for your attention to this matter\") >= 0) {\n    secondPhoneParagraph = p;\n  }\n}\n\nif (!firstPhoneParagraph) {\n  throw new Error(\"Could not find the '{{ requestor.preferred_phone_number }}' paragraph.\");\n}\nif (!secondPhoneParagraph) {\n  throw new Error(\"Could not find the 'Thank you for your attention...' paragraph.\");\n}\n\n// ---- Change 1: requestor.preferred_phone_number -> requestor.phone_number\nlet hits = firstPhoneParagraph.search(\".preferred_phone_number\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one '.preferred_phone_number' in the first paragraph, found \" + hits.items.length);\n}\nhits.items[0].insertText(\".phone_number\", \"Replace\");\nawait context.sync();\n\n// ---- Change 2: drop the stray _GoBack bookmark from its old location\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ---- Change 3: requestor.preferred_phone_number -> requestor. phone_number\nhits = secondPhoneParagraph.search(\".preferred_phone_number\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one '.preferred_phone_number' in the second paragraph, found \" + hits.items.length);\n}\nhits.items[0].insertText(\". phone_number\", \"Replace\");\nawait context.sync();\n\n// Re-insert \"_GoBack\" right after \"requestor.\" (collapsed, zero-width range)\n// in the same placeholder we just rewrote.\nhits = secondPhoneParagraph.search(\"requestor.\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'requestor.' in the second paragraph, found \" + hits.items.length);\n}\nhits.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change summary (per the XML diff):\n#  1. First  \"{{ requestor.preferred_phone_number }}\"  ->  \"{{ requestor.phone_number }}\"\n#  2. Second \"{{ requestor.preferred_phone_number }}\"  ->  \"{{ requestor. phone_number }}\"\n#     (note the stray space after the dot, copied verbatim from the diff)\n#  3. The \"_GoBack\" bookmark is removed from its old spot (right before the\n#     second \"{% if student.grade ... %}\" block) and re-inserted just after\n#     \"requestor.\" inside the second phone-number placeholder above.\n\n$d = $word.ActiveDocument\n\n# Locate the two paragraphs this edit touches, unambiguously, by content.\n$para1 = $null   # \"{{ requestor.preferred_phone_number }} \"\n$para2 = $null   # \"Thank you ... {{ requestor.preferred_phone_number }} with ...\"\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($para1 -eq $null -and $t -like \"*{{ requestor.preferred_phone_number*\") {\n        $para1 = $p\n    }\n    if ($t -like \"*Thank you for your attention to this matter*\") {\n        $para2 = $p\n    }\n}\nif ($para1 -eq $null) {\n    throw \"Could not find the '{{ requestor.preferred_phone_number }}' paragraph.\"\n}\nif ($para2 -eq $null) {\n    throw \"Could not find the 'Thank you for your attention...' paragraph.\"\n}\n\n# ---- Change 1: requestor.preferred_phone_number -> requestor.phone_number\n$r1 = $para1.Range\n$found1 = $r1.Find.Execute(\".preferred_phone_number\", $false, $false, $false, $false, $false, $true, 1, $false, \".phone_number\", 1)\nif (-not $found1) {\n    throw \"Could not find '.preferred_phone_number' in the first paragraph.\"\n}\n\n# ---- Change 2: drop the stray _GoBack bookmark from its old location\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ---- Change 3: requestor.preferred_phone_number -> requestor. phone_number\n$r2 = $para2.Range\n$found2 = $r2.Find.Execute(\".preferred_phone_number\", $false, $false, $false, $false, $false, $true, 1, $false, \". phone_number\", 1)\nif (-not $found2) {\n    throw \"Could not find '.preferred_phone_number' in the second paragraph.\"\n}\n\n# Re-insert \"_GoBack\" right after \"requestor.\" (collapsed, zero-width range)\n# in the same placeholder we just rewrote.\n$r3 = $para2.Range\n$found3 = $r3.Find.Execute(\"requestor.\")\nif (-not $found3) {\n    throw \"Could not find 'requestor.' in the second paragraph.\"\n}\n$r3.Collapse(0)   # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $r3)\n"}
